# Fix syntax error in the "Finalidade" introduction paragraph:
# "...gerenciamento do estoque e geração de relatórios..."
#                            ^^^                      replace " e" with ","
# becomes
# "...gerenciamento do estoque, geração de relatórios..."
#
# The sentence containing this text is split (in the target document) into
# three separate runs:
#   1) "...gerenciamento do estoque"
#   2) ","
#   3) " geração de relatórios de vendas e o gerenciamento de usuários."
# while the run preceding the sentence (", ") must stay untouched.

$d = $word.ActiveDocument

# Locate the unique anchor phrase that starts the sentence we need to edit.
# (The phrase "estoque e geração" also occurs elsewhere in the document, so
# we anchor on the unique "mais especificamente" phrase instead and compute
# offsets relative to it.)
$anchor = $d.Content.Duplicate
$found = $anchor.Find.Execute("mais especificamente", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $sentenceStart = $anchor.Start

    # Offset (within the sentence) of the " e" that must become ",".
    $offset = 140

    $spaceEStart = $sentenceStart + $offset
    $spaceEEnd = $spaceEStart + 2

    # 1) Pre-split the run boundary that precedes the sentence so that the
    #    upcoming text edit does not get merged with the preceding ", " run.
    $boundary = $d.Range($sentenceStart, $sentenceStart + 1)
    $boundary.Bold = 1
    $boundary.Bold = 0

    # 2) Replace " e" with "," (turns " e geração" into ", geração").
    $editRange = $d.Range($spaceEStart, $spaceEEnd)
    $editRange.Text = ","

    # 3) Split the newly inserted "," off into its own run so the text
    #    following it ("  geração...usuários.") remains a distinct run too.
    $commaRange = $d.Range($spaceEStart, $spaceEStart + 1)
    $commaRange.Bold = 1
    $commaRange.Bold = 0
}
